$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (New York -- New York) was previously blank/error-filled because the
# GitHub API rate-limited the scraper; this run succeeded, so fill in the
# scraped values.

# B4: Date Published (serial 44034 == 2020-07-22), same date format used by
# the rest of column B.
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"
$ws.Range("B4").Value = 44034

# C4 / D4: Total Cases / Total Deaths were captured as text in the source
# data. Force text storage, write the digits, then drop the now-superfluous
# formatting so the cell ends up as a plain text value.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "219128"
$ws.Range("C4").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "18803"
$ws.Range("D4").ClearFormats()

# E4-H4: Count/Pct Cases & Deaths Black/AA (numeric).
$ws.Range("E4").Value = 33790
$ws.Range("F4").Value = 5239
$ws.Range("G4").Value = 30.07
$ws.Range("H4").Value = 30.43

# J4: Pct Includes Hispanic Black flips to TRUE.
$ws.Range("J4").Value = $true

# K4 / L4: Count Cases/Deaths Known Race (numeric).
$ws.Range("K4").Value = 112360
$ws.Range("L4").Value = 17217

# O4: status message for the row now reflects a successful fetch.
$ws.Range("O4").Value = "Success!"

# Row 41 (Iowa): Count Cases Black/AA corrected from 3287 to 3288.
$ws.Range("E41").Value = 3288
